{"js": "// The transcription runs for each \"<id>...</id>\" marker were split across\n// multiple w:r runs (e.g. \"<id>\" / \"p061v_1\" / \"</id>\", or even more pieces\n// like \"<id>\" / \"p061\" / \"v\" / \"_2\" / \"</id>\"). The edit collapses each of\n// those spans back into a single run (keeping the formatting of the first,\n// \"<id>\" run) whose text is the full \"<id>...</id>\" string.\nconst ids = [\"p061v_1\", \"p061v_2\", \"p061v_3\", \"p061v_4\"];\n\nconst body = context.document.body;\n\nfor (const id of ids) {\n  const full = `<id>${id}</id>`;\n  const results = body.search(full, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    continue;\n  }\n\n  // Replacing the whole matched range with its own text merges the\n  // multiple underlying runs into a single run that carries the\n  // formatting of the first run in the (former) span.\n  results.items[0].insertText(full, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The transcription runs for each \"<id>...</id>\" marker were split across\n# multiple runs (e.g. \"<id>\" / \"p061v_1\" / \"</id>\", or even more pieces like\n# \"<id>\" / \"p061\" / \"v\" / \"_2\" / \"</id>\"). This script collapses each of\n# those spans back into a single run (keeping the formatting/rsid* of the\n# first, \"<id>\" run) whose text is the full \"<id>...</id>\" string.\n\n$d = $word.ActiveDocument\n\n# Returns a cheap \"signature\" for the character at position $i so we can\n# detect where one run ends and the next begins without needing to touch\n# any internal run collection (the COM object model only exposes Range/\n# Font, not raw runs).\nfunction Get-CharSig($doc, $i) {\n    $c = $doc.Range($i, $i + 1)\n    return ($c.Font.Color.ToString() + \"|\" + $c.Font.Size.ToString() + \"|\" + $c.Font.Name.ToString() + \"|\" + $c.Font.Bold.ToString() + \"|\" + $c.Font.Italic.ToString())\n}\n\nfunction Merge-IdMarker($doc, $idValue) {\n    $needle = \"<id>\" + $idValue + \"</id>\"\n\n    $rng = $doc.Content\n    $found = $rng.Find.Execute($needle)\n    if (-not $found) {\n        return\n    }\n\n    $base = $rng.Start\n    $end = $rng.End\n    $fullText = $rng.Text\n\n    # Walk forward from the start of the match while the per-character\n    # formatting signature stays the same as the very first character \u2014\n    # that span is the existing first run (\"<id>\"), which we leave alone.\n    $sig0 = Get-CharSig $doc $base\n    $firstRunEnd = $base + 1\n    while ($firstRunEnd -lt $end) {\n        $sig = Get-CharSig $doc $firstRunEnd\n        if ($sig -ne $sig0) {\n            break\n        }\n        $firstRunEnd = $firstRunEnd + 1\n    }\n\n    $firstRunLen = $firstRunEnd - $base\n    if ($firstRunLen -ge $fullText.Length) {\n        # Already a single run covering the whole marker - nothing to do.\n        return\n    }\n    $remainder = $fullText.Substring($firstRunLen)\n\n    # Remove the old (separately-formatted / separately-split) runs that\n    # followed the first run, then type the remaining characters right\n    # after the first run so Word appends them into that same run instead\n    # of minting a new one.\n    $tail = $doc.Range($firstRunEnd, $end)\n    $tail.Delete()\n    $insertionPoint = $doc.Range($firstRunEnd, $firstRunEnd)\n    $insertionPoint.InsertAfter($remainder)\n}\n\n$ids = @(\"p061v_1\", \"p061v_2\", \"p061v_3\", \"p061v_4\")\nforeach ($id in $ids) {\n    Merge-IdMarker $d $id\n}\n"}
